$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cl = $m.CustomLayouts.Item(3)
$sh = $cl.Shapes.Item(3)
$sh | Get-Member | Out-String | Write-Host
